$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.892.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.224.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.14%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "85.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.514"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.20%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.471"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.44%  "

$ws.Range("E10").Value = "  +12.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.73%  "

$ws.Range("E13").Value = "  +1.80%  "

$ws.Range("E14").Value = "  +6.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.567.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.221.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.728"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.77%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "39.815.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0888"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.56%  "

$ws.Range("E21").Value = "  +3.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +11.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.80%  "

$ws.Range("E26").Value = "  +5.33%  "

$ws.Range("E27").Value = "  +7.73%  "

$ws.Range("E28").Value = "  +3.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.90%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0716"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +16.81%  "

$ws.Range("E38").Value = "  +3.53%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0995"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.27%  "

$ws.Range("E40").Value = "  +5.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.032.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.61%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +15.51%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0269"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.74%  "

$ws.Range("E48").Value = "  +3.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.439.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.07%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "89.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.20%  "
